$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 190, shifting existing rows 190:205 down to 191:206.
$ws.Rows.Item(190).Insert()

# Populate the new row 190 with the new record (copy constant columns from the
# row below, which used to be row 190 and now sits at row 191).
$ws.Cells.Item(190, 1).Value = 8
$ws.Cells.Item(190, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(190, 3).Value = "Coquimbo"
$ws.Cells.Item(190, 4).Value = 44578
$ws.Cells.Item(190, 5).Value = 4
$ws.Cells.Item(190, 6).Value = 100112003
$ws.Cells.Item(190, 7).Value = "Ajo"
$ws.Cells.Item(190, 8).Value = "Chino"
$ws.Cells.Item(190, 9).Value = "Primera"
$ws.Cells.Item(190, 10).Value = 560
$ws.Cells.Item(190, 11).Value = 18000
$ws.Cells.Item(190, 12).Value = 19000
$ws.Cells.Item(190, 13).Value = 18500
$ws.Cells.Item(190, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(190, 15).Value = "China"
$ws.Cells.Item(190, 16).Value = 1850
$ws.Cells.Item(190, 17).Value = 10
$ws.Cells.Item(190, 18).Value = "Hortaliza"

# Match the date-number style used by column D in the surrounding rows.
$ws.Cells.Item(190, 4).NumberFormat = $ws.Cells.Item(191, 4).NumberFormat
